$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: "ExpectedSourceTemplateFile_Excel" -> "ExpectedSourceTemplateFile"
$ws.Range("I1").Value = "ExpectedSourceTemplateFile"

# --- Clear the old F:H section-parameter cells out of the scenario1 rows (2-4); this data
#     effectively "moves" down into the new scenario2 block below.
$ws.Range("F2:H4").ClearContents()

# --- Row 6: scenario2 header row (mirrors scenario1's row 2) plus the first sub-section + new expected file
$ws.Range("A6").Value = "scenario2"
$ws.Range("B6").Value = "LIVEHTA Automation - Test_NonOncology_Automation_3"
$ws.Range("C6").Value = "LIVEHTA Automation - Test_NonOncology_Automation_3_radio_button"
$ws.Range("D6").Value = "Clinical"
$ws.Range("E6").Value = "Clinical_radio_button"
$ws.Range("F6").Value = "sub_pop_section1"
$ws.Range("G6").Value = "sub_pop_section1_checkbox"
$ws.Range("H6").Value = "sub_pop_section"
$ws.Range("I6").Value = "\Testdata\Non_Oncology\Templates\SLRReport_SourceData\NonOnco_Expected_TestData_with_filters.xlsx"

# --- Row 7: scenario2, second sub-section (intervention)
$ws.Range("A7").Value = "scenario2"
$ws.Range("F7").Value = "intervention_section4"
$ws.Range("G7").Value = "intervention_section4_checkbox"
$ws.Range("H7").Value = "intervention_section"

# --- Row 8: scenario2, third sub-section (study design)
$ws.Range("A8").Value = "scenario2"
$ws.Range("F8").Value = "study_design_section1"
$ws.Range("G8").Value = "study_design_section1_checkbox"
$ws.Range("H8").Value = "study_design_section"

# --- Row 9: scenario2, brand-new fourth sub-section (reported variable)
$ws.Range("A9").Value = "scenario2"
$ws.Range("F9").Value = "reported_variable_section3"
$ws.Range("G9").Value = "reported_variable_section3_checkbox"
$ws.Range("H9").Value = "reported_variable_section"

# --- Re-fit the section-parameter columns now that they hold new/longer text (best-fit,
#     matching the author's saved column widths as closely as this engine's width model allows)
$ws.Columns.Item(6).ColumnWidth = 23.0
$ws.Columns.Item(7).ColumnWidth = 31.833333333333332
$ws.Columns.Item(8).ColumnWidth = 22.0
$ws.Columns.Item(9).ColumnWidth = 91.16666666666667

# --- Move the view/selection the way the author left it
$ws.Range("I1").Select()
